$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# 1. Snapshot the existing data rows (2..12, columns A-E) BEFORE anything is
#    overwritten, because the update re-orders (permutes) these rows and a
#    naive in-place write would clobber source data before it is read.
#    (Column F is intentionally excluded -- it is being removed entirely.)
# --------------------------------------------------------------------------
$old = @{}
for ($r = 2; $r -le 12; $r++) {
    $old[$r] = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2
    )
}

# Remember the style of a date cell in column A so it can be re-applied to
# the brand-new rows appended below (copy/paste-formats keeps the same
# shared cellXfs entry instead of Excel minting a duplicate style).
$ws.Cells.Item(2, 1).Copy()

# --------------------------------------------------------------------------
# 2. Re-write rows 2..12 in their new order.
#    new row -> old row it came from
# --------------------------------------------------------------------------
$rowMap = [ordered]@{
    2  = 10   # 2021-10
    3  = 11   # 2021-11
    4  = 12   # 2021-12
    5  = 2    # 2021-02
    6  = 3    # 2021-03
    7  = 4    # 2021-04
    8  = 5    # 2021-05
    9  = 6    # 2021-06
    10 = 7    # 2021-07
    11 = 8    # 2021-08
    12 = 9    # 2021-09
}

foreach ($newRow in $rowMap.Keys) {
    $vals = $old[$rowMap[$newRow]]

    $ws.Cells.Item($newRow, 1).Value = $vals[0]

    if ("" -eq $vals[1]) {
        $ws.Cells.Item($newRow, 2).Value = ""
    } else {
        $ws.Cells.Item($newRow, 2).Value = $vals[1]
    }

    if ("" -eq $vals[2]) {
        $ws.Cells.Item($newRow, 3).Value = ""
    } else {
        $ws.Cells.Item($newRow, 3).Value = $vals[2]
    }

    $ws.Cells.Item($newRow, 4).Value = $vals[3]
    $ws.Cells.Item($newRow, 5).Value = $vals[4]
}

# --------------------------------------------------------------------------
# 3. Drop column F (智能手表产量) entirely -- header + all data below it.
# --------------------------------------------------------------------------
$ws.Columns.Item(6).Delete()

# --------------------------------------------------------------------------
# 4. Append the new rows (13..28) with fresh data.
# --------------------------------------------------------------------------
$newRows = @(
    @(13, "2022-10", -16.5,  552.6, 5073.2,  -9.3),
    @(14, "2022-11", -42.7,  585,   5659.4,  -14.7),
    @(15, "2022-12", -25,    627,   6277.4,  -16),
    @(16, "2022-02", $null,  $null, 872.9,    5.5),
    @(17, "2022-03", -36.2,  426.2, 1264.2,  -14.8),
    @(18, "2022-04", -2,     501.9, 1766.5,  -11.6),
    @(19, "2022-05", 3.3,    571.9, 2082.2,  -18.7),
    @(20, "2022-06", 1.4,    583.8, 2690.9,  -14.2),
    @(21, "2022-07", -12.9,  487.6, 3179.6,  -14),
    @(22, "2022-08", 16.3,   616.5, 3659.8,  -13.4),
    @(23, "2022-09", 21.4,   695.6, 4355.3,  -9.2),
    @(24, "2023-02", $null,  $null, 697.7,   -31),
    @(25, "2023-03", 16,     590,   1275,    -14),
    @(26, "2023-04", -25.5,  451.4, 1757.8,  -15.3),
    @(27, "2023-05", -24.1,  557,   2313.8,  -8.6),
    @(28, "2023-06", 6.8,    790.8, 3213.7,  -1.7)
)

foreach ($entry in $newRows) {
    $r = $entry[0]

    # Give the date cell the same look (bold, centered, bordered) as the
    # other column-A date cells by pasting the format captured in step 2.
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $entry[1]

    if ($null -eq $entry[2]) {
        $ws.Cells.Item($r, 2).Value = ""
    } else {
        $ws.Cells.Item($r, 2).Value = $entry[2]
    }

    if ($null -eq $entry[3]) {
        $ws.Cells.Item($r, 3).Value = ""
    } else {
        $ws.Cells.Item($r, 3).Value = $entry[3]
    }

    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
}

$excel.CutCopyMode = 0
